# B6-PowerPoint.pptx edit
#
# Source diff does two things:
#   1. Re-points the tableStyleId of the 3 tables (slides 14, 15, 16) from
#      the deck's default "Table_0" style {19738D0E-CAA7-497F-9176-7D5F7BD384D6}
#      to {3D36A673-4D90-4435-A373-20D8AE7A5C1C}.
#   2. Swaps the content of ppt/theme/theme1.xml (the real "Integral" theme
#      used by the slide master) and ppt/theme/theme2.xml (the stock "Office
#      Theme" used only by the notes master) with each other.
#
# (2) can only be reached through the slide master's Theme object in this
# host - there is no COM surface that independently addresses the notes
# master's theme part, so we replicate the reachable half of the swap: push
# the "Office Theme" color values into the slide master's theme color
# scheme (theme1.xml). The theme/clrScheme "name" attributes and
# ppt/theme/theme2.xml itself are not reachable from the object model.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{3D36A673-4D90-4435-A373-20D8AE7A5C1C}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colors (theme1.xml <- "Office Theme" palette) ----------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

$themeColors.Colors(1).RGB  = 0         # dk1      000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388   # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407     # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308  # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink 954F72

Write-Output "done"
